{"js": "// Adds a new bulleted list item \"Sehat walafiat\" after \"Masuk PIMNAS\",\n// matching the formatting (ListParagraph style + same numbering) of the\n// existing list items.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the \"Masuk PIMNAS\" paragraph (the last item in the bullet list).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Masuk PIMNAS\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  // Fallback: use the very last paragraph in the body.\n  target = paragraphs.items[paragraphs.items.length - 1];\n}\n\n// Capture the list this paragraph belongs to, so the new paragraph can be\n// attached to the very same list/numbering instance.\nconst list = target.list;\nlist.load(\"id\");\nawait context.sync();\n\n// Insert a new paragraph right after it, carrying over the same paragraph\n// style and list numbering.\nconst newPara = target.insertParagraph(\"Sehat walafiat\", \"After\");\nnewPara.style = \"List Paragraph\";\nnewPara.attachToList(list.id, 0);\n\nawait context.sync();\n", "ps1": "# Adds a new bulleted list item \"Sehat walafiat\" after \"Masuk PIMNAS\",\n# matching the formatting (ListParagraph style + same numbering) of the\n# existing list items.\n$d = $word.ActiveDocument\n\n# Locate the \"Masuk PIMNAS\" paragraph (the last item of the bullet list).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Masuk PIMNAS\") {\n        $target = $p\n    }\n}\nif ($null -eq $target) {\n    $target = $d.Paragraphs.Last\n}\n\n# Insert a new paragraph right after it; the new paragraph inherits the\n# source paragraph's style/list numbering automatically. Re-fetch the\n# paragraph via the document's collection (the original $target reference\n# becomes stale once the content shifts) and set its text.\n$target.Range.InsertParagraphAfter() | Out-Null\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"Sehat walafiat\"\n"}
